$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36
$ws.Cells.Item(36, 1).Value = 112045414
$ws.Cells.Item(36, 2).Value = 90689
$ws.Cells.Item(36, 3).Value = 'Ovaliderad'
$ws.Cells.Item(36, 4).Value = 'NT'
$ws.Cells.Item(36, 5).Value = 5966
$ws.Cells.Item(36, 6).Value = 'Motaggsvamp'
$ws.Cells.Item(36, 7).Value = 'Sarcodon squamosus'
$ws.Cells.Item(36, 8).Value = '(Schaeff.) Quél.'
$ws.Cells.Item(36, 10).Value = 'fruktkroppar'
$ws.Cells.Item(36, 16).Value = 'Nedre Tetvasseltjärnen (Nedre Tetvasseltjärnen), Dlr'
$ws.Cells.Item(36, 17).Value = 492044.3381435904
$ws.Cells.Item(36, 18).Value = 6785564.065369682
$ws.Cells.Item(36, 19).Value = 15
$ws.Cells.Item(36, 20).Value = 'Dalarna'
$ws.Cells.Item(36, 21).Value = 'Orsa'
$ws.Cells.Item(36, 22).Value = 'Dalarna'
$ws.Cells.Item(36, 23).Value = 'Orsa'
$ws.Cells.Item(36, 25).NumberFormat = "@"
$ws.Cells.Item(36, 25).Value = '2023-09-12'
$ws.Cells.Item(36, 25).Style = "Normal"
$ws.Cells.Item(36, 26).Value = '12:36'
$ws.Cells.Item(36, 27).NumberFormat = "@"
$ws.Cells.Item(36, 27).Value = '2023-09-12'
$ws.Cells.Item(36, 27).Style = "Normal"
$ws.Cells.Item(36, 28).Value = '12:36'
$ws.Cells.Item(36, 30).Value = $false
$ws.Cells.Item(36, 31).Value = $false
$ws.Cells.Item(36, 33).Value = $false
$ws.Cells.Item(36, 49).Value = 'Bo karlstens'
$ws.Cells.Item(36, 50).Value = 'Bo karlstens, Bengt Oldhammer, Janolof Hermansson'

# Row 37
$ws.Cells.Item(37, 1).Value = 112045085
$ws.Cells.Item(37, 2).Value = 88924
$ws.Cells.Item(37, 3).Value = 'Ovaliderad'
$ws.Cells.Item(37, 4).Value = 'LC'
$ws.Cells.Item(37, 5).Value = 256703
$ws.Cells.Item(37, 6).Value = 'Tallfingersvamp'
$ws.Cells.Item(37, 7).Value = 'Ramaria eosanguinea'
$ws.Cells.Item(37, 8).Value = 'R.H.Petersen'
$ws.Cells.Item(37, 9).NumberFormat = "@"
$ws.Cells.Item(37, 9).Value = '3'
$ws.Cells.Item(37, 9).Style = "Normal"
$ws.Cells.Item(37, 16).Value = 'Nedre Tetvasseltjärnen (Nedre Tetvasseltjärnen), Dlr'
$ws.Cells.Item(37, 17).Value = 492013.8535130407
$ws.Cells.Item(37, 18).Value = 6785532.323436439
$ws.Cells.Item(37, 19).Value = 5
$ws.Cells.Item(37, 20).Value = 'Dalarna'
$ws.Cells.Item(37, 21).Value = 'Orsa'
$ws.Cells.Item(37, 22).Value = 'Dalarna'
$ws.Cells.Item(37, 23).Value = 'Orsa'
$ws.Cells.Item(37, 25).NumberFormat = "@"
$ws.Cells.Item(37, 25).Value = '2023-09-12'
$ws.Cells.Item(37, 25).Style = "Normal"
$ws.Cells.Item(37, 26).Value = '12:14'
$ws.Cells.Item(37, 27).NumberFormat = "@"
$ws.Cells.Item(37, 27).Value = '2023-09-12'
$ws.Cells.Item(37, 27).Style = "Normal"
$ws.Cells.Item(37, 28).Value = '12:14'
$ws.Cells.Item(37, 30).Value = $false
$ws.Cells.Item(37, 31).Value = $false
$ws.Cells.Item(37, 33).Value = $false
$ws.Cells.Item(37, 49).Value = 'Bo karlstens'
$ws.Cells.Item(37, 50).Value = 'Bo karlstens, Bengt Oldhammer, Janolof Hermansson'

# Row 38
$ws.Cells.Item(38, 1).Value = 112045302
$ws.Cells.Item(38, 2).Value = 90658
$ws.Cells.Item(38, 3).Value = 'Ovaliderad'
$ws.Cells.Item(38, 4).Value = 'NT'
$ws.Cells.Item(38, 5).Value = 4361
$ws.Cells.Item(38, 6).Value = 'Orange taggsvamp'
$ws.Cells.Item(38, 7).Value = 'Hydnellum aurantiacum'
$ws.Cells.Item(38, 8).Value = '(Batsch:Fr.) P.Karst.'
$ws.Cells.Item(38, 16).Value = 'Nedre Tetvasseltjärnen (Nedre Tetvasseltjärnen), Dlr'
$ws.Cells.Item(38, 17).Value = 492044.3381435904
$ws.Cells.Item(38, 18).Value = 6785564.065369682
$ws.Cells.Item(38, 19).Value = 5
$ws.Cells.Item(38, 20).Value = 'Dalarna'
$ws.Cells.Item(38, 21).Value = 'Orsa'
$ws.Cells.Item(38, 22).Value = 'Dalarna'
$ws.Cells.Item(38, 23).Value = 'Orsa'
$ws.Cells.Item(38, 25).NumberFormat = "@"
$ws.Cells.Item(38, 25).Value = '2023-09-12'
$ws.Cells.Item(38, 25).Style = "Normal"
$ws.Cells.Item(38, 26).Value = '12:27'
$ws.Cells.Item(38, 27).NumberFormat = "@"
$ws.Cells.Item(38, 27).Value = '2023-09-12'
$ws.Cells.Item(38, 27).Style = "Normal"
$ws.Cells.Item(38, 28).Value = '12:27'
$ws.Cells.Item(38, 30).Value = $false
$ws.Cells.Item(38, 31).Value = $false
$ws.Cells.Item(38, 33).Value = $false
$ws.Cells.Item(38, 49).Value = 'Bo karlstens'
$ws.Cells.Item(38, 50).Value = 'Bo karlstens, Bengt Oldhammer, Janolof Hermansson'

# Row 39
$ws.Cells.Item(39, 1).Value = 112045406
$ws.Cells.Item(39, 2).Value = 90682
$ws.Cells.Item(39, 3).Value = 'Ovaliderad'
$ws.Cells.Item(39, 4).Value = 'NT'
$ws.Cells.Item(39, 5).Value = 2059
$ws.Cells.Item(39, 6).Value = 'Skrovlig taggsvamp'
$ws.Cells.Item(39, 7).Value = 'Hydnellum scabrosum'
$ws.Cells.Item(39, 8).Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Cells.Item(39, 16).Value = 'Nedre Tetvasseltjärnen (Nedre Tetvasseltjärnen), Dlr'
$ws.Cells.Item(39, 17).Value = 492044.3381435904
$ws.Cells.Item(39, 18).Value = 6785564.065369682
$ws.Cells.Item(39, 19).Value = 15
$ws.Cells.Item(39, 20).Value = 'Dalarna'
$ws.Cells.Item(39, 21).Value = 'Orsa'
$ws.Cells.Item(39, 22).Value = 'Dalarna'
$ws.Cells.Item(39, 23).Value = 'Orsa'
$ws.Cells.Item(39, 25).NumberFormat = "@"
$ws.Cells.Item(39, 25).Value = '2023-09-12'
$ws.Cells.Item(39, 25).Style = "Normal"
$ws.Cells.Item(39, 26).Value = '12:36'
$ws.Cells.Item(39, 27).NumberFormat = "@"
$ws.Cells.Item(39, 27).Value = '2023-09-12'
$ws.Cells.Item(39, 27).Style = "Normal"
$ws.Cells.Item(39, 28).Value = '12:36'
$ws.Cells.Item(39, 30).Value = $false
$ws.Cells.Item(39, 31).Value = $false
$ws.Cells.Item(39, 33).Value = $false
$ws.Cells.Item(39, 49).Value = 'Bo karlstens'
$ws.Cells.Item(39, 50).Value = 'Bo karlstens, Bengt Oldhammer, Janolof Hermansson'
